$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "CDTSFILIA ; "
$ws.Range("B13").Value = "CCODDFT ; `nCXLAPTZ ; `nCYLAPTZ"
$ws.Range("B18").Value = "CDTSCVT ; `nCDTMCVT"
$ws.Range("A29").Value = "ILTASIT ; `nILTASEU"
$ws.Range("C34").Value = "TA_SEUIL ; `nTA_INFOS_SEUIL"
$ws.Range("D34").Value = "DATE_SAISIE ; `nDATE_SAISIE"

$ws.Range("E20").Select()
